$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-30 Monday" "2024-10-01 Tuesday"

Replace-Text "787×7=" "109×4="
Replace-Text "346×2=" "553×3="
Replace-Text "467×4=" "453×9="
Replace-Text "780×7=" "123×4="
Replace-Text "208×6=" "730×3="
Replace-Text "910×6=" "507×3="
Replace-Text "317×9=" "394×2="
Replace-Text "547×7=" "895×5="
Replace-Text "652×2=" "899×6="
Replace-Text "406×5=" "680×3="
Replace-Text "491×8=" "633×4="
Replace-Text "985×4=" "599×6="
Replace-Text "356×8=" "784×6="
Replace-Text "435×4=" "347×3="
Replace-Text "447×9=" "232×6="
Replace-Text "506×8=" "229×4="
Replace-Text "354×6=" "337×2="
Replace-Text "977×8=" "592×9="
Replace-Text "411×6=" "816×5="
Replace-Text "658×8=" "721×9="
Replace-Text "869×6=" "533×4="
Replace-Text "730×7=" "563×2="
Replace-Text "506×7=" "498×8="
Replace-Text "497×8=" "148×8="
Replace-Text "841×3=" "130×2="
